# Applies the "removed some spelling and grammar errors with the auto checker"
# edit: adjusts estimated/actual hours in the Sprint 4 and Sprint 5 planning
# tables and updates the corresponding explanatory remarks in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 4 table (rows 13-19) ---

# Row 14: "Feedback verwerken." - actual hours spent (F) changes from 0 to 15
$ws.Range("F14").Value2 = 15

# Row 17: "Container opdracht verbeteren." - estimated hours (E) changes 40 -> 50
$ws.Range("E17").Value2 = 50

# Row 19: "Documentatie bijwerken." - priority (D) changes 3 -> 2,
# estimated hours (E) changes 4 -> 5
$ws.Range("D19").Value2 = 2
$ws.Range("E19").Value2 = 5

# --- Sprint 5 table (rows 25-30) ---

# Row 26: "Feedback verwerken." - estimated hours (E) changes 30 -> 40
$ws.Range("E26").Value2 = 40

# Row 27: "Fouthandeling toevoegen." - estimated hours (E) changes 8 -> 10
$ws.Range("E27").Value2 = 10

# Row 28: "SOLID" - remark updated from a maximum of 10 uur to 15 uur.
$ws.Range("H28").Value2 = "Dit kan erg lang doorblijven gaan dus ik zet er een maximum van 15 uur voor deze sprint op."

# Row 15: "Logica toevoegen voor het toevoegen van films, halls, screenings,
# users, tasks en employees." - estimated hours (E) changes 15 -> 10, and the
# remark explaining the estimate is updated from 20 uur to 10 uur.
$ws.Range("E15").Value2 = 10
$ws.Range("H15").Value2 = "Dit is een redelijk breed onderwerp en ik heb nog niet precies bepaald aan welke onderdelen ik tijd wil gaan besteden. Daarom heb ik 10 uur hiervoor vrijgemaakt."

# Row 16: "Logica toevoegen voor het automatisch aanmaken van een tasktype..."
# estimated hours (E) changes 10 -> 5, remark updated from 10 uur to 5 uur.
$ws.Range("E16").Value2 = 5
$ws.Range("H16").Value2 = "Dit zal niet extreem lastig worden maar voor het geval het moeilijker blijkt te zijn heb ik toch 5 uur vrij gemaakt."

# --- Restore the last active selection recorded in the saved workbook ---
$ws.Range("H18").Select()
